$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Range("H96").Value = 1353.25
$ws.Range("I96").Value = 550
$ws.Range("J96").Value = 1513.9
$ws.Range("K96").Value = 1650
$ws.Range("L96").Value = 4541.700000000001
$ws.Range("M96").Value = -277
$ws.Range("N96").Value = -7287.700000000001

# Row 104
$ws.Range("H104").Value = 373
$ws.Range("I104").Value = 86.5
$ws.Range("J104").Value = 946
$ws.Range("K104").Value = 259.5
$ws.Range("L104").Value = 2838
$ws.Range("M104").Value = 1487.5

# Row 125
$ws.Range("H125").Value = 2480921.2
$ws.Range("I125").Value = 523
$ws.Range("J125").Value = 3721120.2
$ws.Range("K125").Value = 4707
$ws.Range("L125").Value = 33490081.8
$ws.Range("M125").Value = -2247
$ws.Range("N125").Value = -33495001.8

# Row 141
$ws.Range("H141").Value = 1869.1034
$ws.Range("I141").Value = 1426.8846
$ws.Range("J141").Value = 5701.6665
$ws.Range("K141").Value = 4280.6538
$ws.Range("L141").Value = 17104.9995
$ws.Range("M141").Value = 899.3462

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2009168.8
$ws.Range("I32").Value = 2280174.8
$ws.Range("J32").Value = 21791.25
$ws.Range("K32").Value = 2280174.8
$ws.Range("L32").Value = 21791.25
$ws.Range("M32").Value = -2279887.8
$ws.Range("N32").Value = -22365.25

# Row 37
$ws.Range("H37").Value = 10064
$ws.Range("I37").Value = 5444.4443
$ws.Range("J37").Value = 14683.556
$ws.Range("K37").Value = 5444.4443
$ws.Range("L37").Value = 14683.556
$ws.Range("M37").Value = -5171.4443
$ws.Range("N37").Value = -15229.556

# Row 74
$ws.Range("H74").Value = 29488864
$ws.Range("I74").Value = 25642072
$ws.Range("J74").Value = 41029240
$ws.Range("K74").Value = 25642072
$ws.Range("L74").Value = 41029240
$ws.Range("M74").Value = -25641198

# Row 77
$ws.Range("H77").Value = 29488864
$ws.Range("I77").Value = 25642072
$ws.Range("J77").Value = 41029240
$ws.Range("K77").Value = 128210360
$ws.Range("L77").Value = 205146200
$ws.Range("M77").Value = -128205992

# Row 97
$ws.Range("H97").Value = 514.64514
$ws.Range("I97").Value = 352.3
$ws.Range("J97").Value = 809.8182
$ws.Range("K97").Value = 352.3
$ws.Range("L97").Value = 809.8182
$ws.Range("M97").Value = 143.7
$ws.Range("N97").Value = -1801.8182

# Row 105
$ws.Range("N105").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0

# Row 132
$ws.Range("H132").Value = 17663204
$ws.Range("I132").Value = 20437446
$ws.Range("J132").Value = 5053006.5
$ws.Range("K132").Value = 61312338
$ws.Range("L132").Value = 15159019.5
$ws.Range("M132").Value = -61309808
$ws.Range("N132").Value = -15164079.5

$ws = $wb.Worksheets.Item("BSM")
# Row 69
$ws.Range("M69").ClearContents()
$ws.Range("H69").Value = 24000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 24000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 24000
$ws.Range("N69").Value = -25622

# Row 72
$ws.Range("M72").ClearContents()
$ws.Range("H72").Value = 24000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 24000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 72000
$ws.Range("N72").Value = -80112

# Row 134
$ws.Range("H134").Value = 10389072
$ws.Range("I134").Value = 14414272
$ws.Range("J134").Value = 38558.5
$ws.Range("K134").Value = 43242816
$ws.Range("L134").Value = 115675.5
$ws.Range("M134").Value = -43240281
$ws.Range("N134").Value = -120745.5

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 402440.34
$ws.Range("I2").Value = 3660.5
$ws.Range("J2").Value = 1200000
$ws.Range("K2").Value = 3660.5
$ws.Range("L2").Value = 1200000
$ws.Range("M2").Value = -3547.5
$ws.Range("N2").Value = -1200226

# Row 58
$ws.Range("H58").Value = 930648
$ws.Range("I58").Value = 3435.8647
$ws.Range("J58").Value = 3789552
$ws.Range("K58").Value = 3435.8647
$ws.Range("L58").Value = 3789552
$ws.Range("M58").Value = -3232.8647
$ws.Range("N58").Value = -3789958

# Row 60
$ws.Range("H60").Value = 8201.286
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 8201.286
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 8201.286
$ws.Range("N60").Value = -9223.286

# Row 68
$ws.Range("H68").Value = 18149.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 18149.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18149.75
$ws.Range("N68").Value = -19647.75

# Row 71
$ws.Range("H71").Value = 18149.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 18149.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 54449.25
$ws.Range("N71").Value = -61937.25

# Row 132
$ws.Range("H132").Value = 1455.8276
$ws.Range("I132").Value = 1062.238
$ws.Range("J132").Value = 2489
$ws.Range("K132").Value = 3186.714
$ws.Range("L132").Value = 7467
$ws.Range("M132").Value = -656.7139999999999
$ws.Range("N132").Value = -12527

# Row 134
$ws.Range("H134").Value = 1113511.1
$ws.Range("I134").Value = 1516.4231
$ws.Range("J134").Value = 4004697.2
$ws.Range("K134").Value = 4549.2693
$ws.Range("L134").Value = 12014091.6
$ws.Range("M134").Value = -2014.2693

# Row 136
$ws.Range("H136").Value = 930648
$ws.Range("I136").Value = 3435.8647
$ws.Range("J136").Value = 3789552
$ws.Range("K136").Value = 10307.5941
$ws.Range("L136").Value = 11368656
$ws.Range("M136").Value = -7757.5941
$ws.Range("N136").Value = -11373756

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 497.0465
$ws.Range("I122").Value = 292.14285
$ws.Range("J122").Value = 1393.5
$ws.Range("K122").Value = 2629.28565
$ws.Range("L122").Value = 12541.5
$ws.Range("M122").Value = -179.2856500000003
$ws.Range("N122").Value = -17441.5

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 307.21054
$ws.Range("I107").Value = 186.91667
$ws.Range("J107").Value = 513.4286
$ws.Range("K107").Value = 186.91667
$ws.Range("L107").Value = 513.4286
$ws.Range("M107").Value = 1733.08333
$ws.Range("N107").Value = -4353.4286

# Row 132
$ws.Range("H132").Value = 6511833
$ws.Range("I132").Value = 6517526
$ws.Range("J132").Value = 6496380
$ws.Range("K132").Value = 19552578
$ws.Range("L132").Value = 19489140
$ws.Range("M132").Value = -19550048
$ws.Range("N132").Value = -19494200

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4203939.5
$ws.Range("I132").Value = 7144272
$ws.Range("J132").Value = 3464.2144
$ws.Range("K132").Value = 21432816
$ws.Range("L132").Value = 10392.6432
$ws.Range("M132").Value = -21430286
$ws.Range("N132").Value = -15452.6432

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 305.5
$ws.Range("I113").Value = 233.26666
$ws.Range("J113").Value = 404
$ws.Range("K113").Value = 699.79998
$ws.Range("L113").Value = 1212
$ws.Range("M113").Value = 1470.20002
$ws.Range("N113").Value = -5552

# Row 132
$ws.Range("H132").Value = 995823.0600000001
$ws.Range("I132").Value = 4739.091
$ws.Range("J132").Value = 1834432.6
$ws.Range("K132").Value = 14217.273
$ws.Range("L132").Value = 5503297.800000001
$ws.Range("M132").Value = -11687.273
$ws.Range("N132").Value = -5508357.800000001

# Row 136
$ws.Range("H136").Value = 1264.902
$ws.Range("I136").Value = 875.125
$ws.Range("J136").Value = 2682.2727
$ws.Range("K136").Value = 2625.375
$ws.Range("L136").Value = 8046.8181
$ws.Range("M136").Value = -75.375
$ws.Range("N136").Value = -13146.8181
